$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "price_poachers_per_totoaba"
$ws.Range("F2:F7").Value = 16045
$ws.Range("G1").Select()
